$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (all values must remain plain text
# to match the original inlineStr / text cell formatting).
$updates = @{
    "D2" = "324.89"
    "E2" = "-1.45%"
    "E3" = "-1.80%"
    "D4" = "5.690"
    "E4" = "7.66%"
    "D5" = "0.08002"
    "E5" = "-1.29%"
    "D6" = "2.012"
    "E6" = "4.66%"
    "D7" = "4.500"
    "E7" = "-0.35%"
    "D8" = "8.626"
    "E8" = "-0.17%"
    "D10" = "0.9248"
    "D11" = "0.1231"
    "E11" = "-8.53%"
    "D12" = "0.1972"
    "E12" = "0.15%"
    "D13" = "8.699"
    "E13" = "21.45%"
    "D14" = "0.09228"
    "E14" = "-0.16%"
    "D15" = "0.03626"
    "E15" = "1.57%"
    "E16" = "9.41%"
    "D17" = "0.001300"
    "E17" = "-2.24%"
    "D18" = "0.006132"
    "E18" = "-2.04%"
    "E19" = "-0.71%"
    "E20" = "-1.38%"
    "D21" = "0.1372"
    "E21" = "3.77%"
    "D22" = "0.2413"
    "E22" = "-5.72%"
    "D23" = "0.04397"
    "E23" = "-0.78%"
    "D24" = "0.001262"
    "E24" = "3.45%"
    "D25" = "0.004629"
    "E25" = "6.30%"
    "E26" = "-3.28%"
    "D39" = "0.02487"
    "E39" = "-0.65%"
    "D40" = "0.05337"
    "E40" = "2.82%"
    "D41" = "0.007458"
    "E41" = "-2.68%"
    "D42" = "0.009708"
    "E42" = "6.78%"
    "D43" = "0.1407"
    "E43" = "-1.42%"
    "E44" = "-2.42%"
    "D45" = "0.01115"
    "E45" = "9.80%"
    "D46" = "0.00006757"
    "E46" = "1.63%"
    "E47" = "0.09%"
    "E48" = "-11.15%"
    "D49" = "0.002292"
    "E49" = "-4.49%"
    "D50" = "0.00002101"
    "E50" = "0.09%"
    "D51" = "0.0002001"
    "E51" = "0.09%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
